$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string into a
# floating point number (losing exact text representation / trailing zeros),
# which would not match the original inline-string text cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "42.948.44"
$ws.Range("D3").Value = "2.295.79"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "299.67"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "97.57"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").Value = "0.516"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "35.78"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "17.70"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "6.78"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "2.654.80"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "2.295.09"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "0.778"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "42.871.50"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "12.64"
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "6.10"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").Value = "68.08"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "241.45"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "25.12"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "166.44"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "9.06"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "32.95"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "17.60"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "0.0687"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "2.004.89"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "17.29"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "53.52"
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.523.35"
$ws.Range("E51").Value = "  -1.00%  "
